$d = $word.ActiveDocument

# Locate the paragraph that starts with "Iconos, imágenes y logotipo:" - the
# very next paragraph is "Para este contenido usare ...". Both need to be
# merged into a single paragraph (the text was split across two paragraphs
# by mistake), matching the corrected reading order.
$startPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Iconos, im*genes y logotipo:*") {
        $startPara = $p
        break
    }
}

if ($startPara -eq $null) {
    throw "Could not locate the 'Iconos, imágenes y logotipo:' paragraph"
}

$endPara = $startPara.Next()
if ($endPara -eq $null -or -not ($endPara.Range.Text -like "Para este contenido usare p*ginas*")) {
    throw "Could not locate the 'Para este contenido usare...' paragraph"
}

$mergedRange = $d.Range($startPara.Range.Start, $endPara.Range.End)

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="10F23E9B" w14:textId="38739C0D" w:rsidR="00DC63FE" w:rsidRDefault="00DC63FE" w:rsidP="00DC63FE" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Iconos, imágenes y logotipo</w:t></w:r>' + `
    '<w:r><w:t>:</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> p</w:t></w:r>' + `
    '<w:r><w:t>ara este contenido usare páginas que dispongan de imágenes que puedan ser usadas libremente como por ejemplo unsplash</w:t></w:r>' + `
    '<w:r w:rsidR="008A411B"><w:t xml:space="preserve"> para fotografías y para iconos/logotipo undraw y boxicons, estas páginas son de gran ayuda a la hora de usar este tipo de elementos y tienen una fácil accesibilidad.</w:t></w:r>' + `
    '</w:p>'

[void]$mergedRange.InsertXML($newXml)
